$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 1.68
$ws.Range("H2").Value = 2.46
$ws.Range("J2").Value = 2.5
$ws.Range("W2").Value = 2.46

# Row 3 updates
$ws.Range("F3").Value = 1.7
$ws.Range("G3").Value = 1.86
$ws.Range("H3").Value = 4.6
$ws.Range("I3").Value = 5.5
$ws.Range("J3").Value = 3.8
$ws.Range("K3").Value = 5.1
$ws.Range("P3").Value = 2.08
$ws.Range("Q3").Value = 1.63

# Row 8 updates
$ws.Range("F8").Value = 2.2
$ws.Range("G8").Value = 2.4
$ws.Range("H8").Value = 3.85
$ws.Range("I8").Value = 4.3
$ws.Range("P8").Value = 1.55
$ws.Range("Q8").Value = 2.5
